{"js": "// The author repositioned the cursor's last-edit marker (\"_GoBack\") from the\n// end of the \"Connecting Cables\" intro paragraph to the \"* Symbol\n// Multiplication\" heading (splitting \"* \" and \"Symbol Multiplication\" into\n// separate runs around the bookmark). Word keeps \"_GoBack\" as a singleton\n// bookmark, so re-adding it at the new location implicitly removes the old\n// one; we do that explicitly here since the Office.js shim does not merge\n// same-name bookmarks automatically.\n\n// Remove any existing \"_GoBack\" bookmark (the stale one after \"...indexed\n// from 1 to n. \" in the Connecting Cables section).\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// Find the heading text \"Symbol Multiplication\" (part of \"* Symbol\n// Multiplication\") and wrap it with a fresh \"_GoBack\" bookmark.\nconst results = context.document.body.search(\"Symbol Multiplication\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# The author repositioned the cursor's last-edit marker (\"_GoBack\") from the\n# end of the \"Connecting Cables\" intro paragraph to the \"* Symbol\n# Multiplication\" heading (splitting \"* \" and \"Symbol Multiplication\" into\n# separate runs around the bookmark). \"_GoBack\" is a singleton bookmark in\n# Word, so Bookmarks.Add with that name automatically replaces whichever\n# range previously held it.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"Symbol Multiplication\"\n$find.MatchCase = $true\n$find.Execute() | Out-Null\n\nif ($find.Found) {\n    $rng = $find.Parent.Duplicate\n    $d.Bookmarks.Add(\"_GoBack\", $rng) | Out-Null\n}\n"}
